# Update instance_8_times.xlsx data after modifying penalties / fixing
# assignments and new population generation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-11 (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$ws.Range("A2").Value = 2

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 10

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 15

$ws.Range("A7").Value = 8
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 20

$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5

$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 11

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5

$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 10

# New rows 12 and 13
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 15

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 21
